$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.643.51'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '1.632.53'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.97'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.499'
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.15'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("E11").Value = '  +3.56%  '
$ws.Range("D12").Value = '1.860.91'
$ws.Range("D13").Value = '1.641.73'
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  +1.87%  '
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '26.645.36'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.23'
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.93'
$ws.Range("E19").Value = '  +8.08%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.39'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  +2.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.03'
$ws.Range("E25").Value = '  +2.28%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.88'
$ws.Range("E28").Value = '  +4.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").Value = '1.209.49'
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("E37").Value = '  +5.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.809'
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.501'
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.27'
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.41'
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.791'
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '1.773.69'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.72'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.55'
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.67'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0514'
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.63'
$ws.Range("E49").Value = '  +4.74%  '
$ws.Range("E50").Value = '  +0.32%  '
